$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Content")

# Insert 3 new rows before current row 3 ("End"), pushing the old row 3 ("End") down to row 6,
# and shifting the rows below it (13-16) down to (16-19).
$ws.Rows("3:5").Insert()

# New rows 3-5: additional "Prod" content rows (Environment/Brand/Campaign)
$ws.Range("A3").Value = "Prod"
$ws.Range("B3").Value = "MeaningfulBeauty"
$ws.Range("C3").Value = "Core"

$ws.Range("A4").Value = "Prod"
$ws.Range("B4").Value = "WestmoreBeauty"
$ws.Range("C4").Value = "Core"

$ws.Range("A5").Value = "Prod"
$ws.Range("B5").Value = "JLoBeauty"
$ws.Range("C5").Value = "Core"

# Row 6 already holds "End" (shifted down from old row 3) - leave as is.

# Append a new block of rows (15-19) repeating the same content pattern.
# Row 15 sits in a previously-empty area, so give it the shaded row formatting
# used by the rest of the data rows (copy format from row 2).
$ws.Range("A2:C2").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A15").Value = "Prod"
$ws.Range("B15").Value = "CrepeErase"
$ws.Range("C15").Value = "Core"

$ws.Range("A16").Value = "Prod"
$ws.Range("B16").Value = "MeaningfulBeauty"
$ws.Range("C16").Value = "core_full_30_90"

$ws.Range("A17").Value = "Prod"
$ws.Range("B17").Value = "WestmoreBeauty"
$ws.Range("C17").Value = "Core"

$ws.Range("A18").Value = "Prod"
$ws.Range("B18").Value = "JLoBeauty"
$ws.Range("C18").Value = "Core"

$ws.Range("A19").Value = "End"

$ws.Range("A6:XFD6").Select()
